# Update the "Survey Data" sheet header:
#  - Row 1 becomes a set of group headers ("Personal Information" / "Rankings" /
#    "Product Features") spanning columns A:BM, with the non-label cells left
#    blank (they are merged/spanned visually, but kept as separate blank cells
#    here so the used range extends out to column BM).
#  - Row 2 becomes the real column headers: Name / Phone / Bike Type, followed
#    by 15 ranking columns ("Row 1" .. "Row 15").
#  - The old sample data rows (previously rows 2-6) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start from a clean slate: wipe every existing cell (A1:E6) ------------
$ws.Range("A1:E6").ClearContents()

# --- Row 1: top-level group headers -----------------------------------------
$ws.Range("A1").Value = "Personal Information"
$ws.Range("D1").Value = "Rankings"
$ws.Range("T1").Value = "Product Features"

# The remaining header cells on row 1 (B1:C1, E1:S1, U1:BM1) are blank, but we
# still want them to be "real" cells so the sheet's used range extends to
# column BM (matching the new header band spanned across A:BM). Touching a
# formatting property (here: an explicit "no border", which is already the
# default) registers the cell in the sheet without writing any visible text
# or introducing a new cell style.
$ws.Range("B1:C1").Borders.LineStyle = -4142
$ws.Range("E1:S1").Borders.LineStyle = -4142
$ws.Range("U1:BM1").Borders.LineStyle = -4142

# --- Row 2: real column headers ---------------------------------------------
$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "Phone"
$ws.Range("C2").Value = "Bike Type"

for ($i = 1; $i -le 15; $i++) {
    $col = $ws.Cells.Item(2, 3 + $i)
    $col.Value = "Row $i"
}

Write-Host "done"
